# Clean up embedded newlines (and one stray trailing asterisk) inside
# several cell values on the vaccine price-list sheet.  Excel will
# naturally dedupe the shared-string table against the now-identical
# single-line strings when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "DTaP-Hep B-IPV"
$ws.Range("A7").Value = "DTaP-Hep B-IPV"

$ws.Range("D8").Value = "5 x 1 dose vials"
$ws.Range("H8").Value = "Aventis Pateur"

$ws.Range("D11").Value = "10 x 1 dose vials"
$ws.Range("D12").Value = "10 x 1 dose vial"

$ws.Range("D18").Value = "5 x 1 dose vials"

$ws.Range("D24").Value = "10 x 1dose vial 5 x 1 dose TIP-LOK"

$ws.Range("D38").Value = "10 x 1 dose vials"
$ws.Range("D39").Value = "5 x 1 dose vial"
$ws.Range("D40").Value = "5 x 1 dose vial"
$ws.Range("H40").Value = "Aventis Pasteur"

$ws.Range("D41").Value = "10 dose vials"
$ws.Range("H41").Value = "Aventis Pasteur"

$ws.Range("D42").Value = "10 dose vial"
